$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows (10-12), matching the existing table's layout:
# A=Principle, B=Start Principle, C=BuyPrice, D=SellPrice,
# E=IsShortSell (bool), F=Price Change %, G=Date (date-formatted), H=Profitable (bool)

$rows = @(
    @{ Row = 10; A = 9581.06;            B = 9597.3799999999992; C = 78.05;            D = 77.92; E = $false; F = -0.17; G = 42613.765462962961; H = $false },
    @{ Row = 11; A = 9659.6200000000008; B = 9581.06;            C = 77.739999999999995; D = 78.38; E = $false; F = 0.82;  G = 42614.672638888886; H = $true  },
    @{ Row = 12; A = 9543.7000000000007; B = 9659.6200000000008; C = 78.36;            D = 77.42; E = $false; F = -1.2;  G = 42615.750034722223; H = $false }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = $r.G
    $gCell.NumberFormat = "m/d/yy h:mm"

    $ws.Cells.Item($row, 8).Value = $r.H
}
